$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8166584968566895
$ws.Range("B1").Value = 3.022672653198242
$ws.Range("C1").Value = 3.02791166305542
$ws.Range("D1").Value = 2.568146228790283
$ws.Range("E1").Value = 2.185593843460083
